# Generate Report for Archive
#
# The localization status for the sample file moved from
# "Ready for handoff" to "In Translation". That status string shows up
# on the Overview sheet (zh-cn / de-de status columns, E2 & F2) as well
# as on each per-locale sheet's "Status" column (C2 on "zh-cn" and
# "de-de"). Update every occurrence so they keep referring to one
# shared string.
$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus

# The status columns got noticeably narrower (from ~17.22 chars to
# ~13.41 chars) once the new, shorter status text was applied - mirror
# that resize on every sheet that carries a status column.
$narrowWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $narrowWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $narrowWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $narrowWidth
